$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -5
$ws.Range("F11").Value = 2
